$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = "'64.510.75"
$c.Style = 'Normal'
$ws.Range('E2').Value = '  +3.11%  '
$c = $ws.Range('D3')
$c.Value = "'3.459.38"
$c.Style = 'Normal'
$ws.Range('E3').Value = '  +3.77%  '
$ws.Range('E4').Value = '  +0.02%  '
$c = $ws.Range('D5')
$c.Value = "'577.21"
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +4.36%  '
$c = $ws.Range('D6')
$c.Value = "'158.19"
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +4.10%  '
$c = $ws.Range('D7')
$c.Value = "'0.999"
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +0.05%  '
$c = $ws.Range('D8')
$c.Value = "'3.459.54"
$c.Style = 'Normal'
$ws.Range('E8').Value = '  +3.69%  '
$ws.Range('E9').Value = '  +5.01%  '
$c = $ws.Range('D10')
$c.Value = "'7.55"
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +0.70%  '
$ws.Range('E11').Value = '  +6.42%  '
$c = $ws.Range('D12')
$c.Value = "'0.444"
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +2.25%  '
$c = $ws.Range('D13')
$c.Value = "'4.057.82"
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +3.88%  '
$c = $ws.Range('D14')
$c.Value = "'0.135"
$c.Style = 'Normal'
$ws.Range('E14').Value = '  -2.02%  '
$ws.Range('E15').Value = '  +9.98%  '
$c = $ws.Range('D16')
$c.Value = "'27.79"
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +3.66%  '
$c = $ws.Range('D17')
$c.Value = "'64.540.93"
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +3.11%  '
$c = $ws.Range('D18')
$c.Value = "'3.455.50"
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +2.86%  '
$ws.Range('E19').Value = '  -1.12%  '
$c = $ws.Range('D20')
$c.Value = "'14.36"
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +4.99%  '
$c = $ws.Range('D21')
$c.Value = "'395.18"
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +2.43%  '
$c = $ws.Range('D22')
$c.Value = "'8.49"
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +0.72%  '
$ws.Range('E23').Value = '  +1.49%  '
$c = $ws.Range('D24')
$c.Value = "'72.89"
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +3.17%  '
$c = $ws.Range('D25')
$c.Value = "'0.999"
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -0.27%  '
$ws.Range('E26').Value = '  +28.82%  '
$c = $ws.Range('D27')
$c.Value = "'9.64"
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +9.05%  '
$ws.Range('E28').Value = '  +1.67%  '
$ws.Range('E29').Value = '  -0.14%  '
$c = $ws.Range('D30')
$c.Value = "'6.13"
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +10.57%  '
$ws.Range('E31').Value = '  +7.56%  '
$c = $ws.Range('D32')
$c.Value = "'6.68"
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +4.53%  '
$c = $ws.Range('D33')
$c.Value = "'2.03"
$c.Style = 'Normal'
$ws.Range('E33').Value = '  +2.97%  '
$c = $ws.Range('D34')
$c.Value = "'23.77"
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +3.50%  '
$ws.Range('E35').Value = '  -0.05%  '
$c = $ws.Range('D36')
$c.Value = "'7.04"
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +5.36%  '
$c = $ws.Range('D37')
$c.Value = "'160.65"
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -0.29%  '
$ws.Range('E38').Value = '  -0.35%  '
$c = $ws.Range('D39')
$c.Value = "'0.0785"
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +6.42%  '
$c = $ws.Range('D40')
$c.Value = "'1.87"
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +0.39%  '
$c = $ws.Range('D41')
$c.Value = "'27.58"
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +0.18%  '
$c = $ws.Range('D42')
$c.Value = "'2.938.05"
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +2.93%  '
$c = $ws.Range('D43')
$c.Value = "'0.0323"
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +3.06%  '
$c = $ws.Range('D44')
$c.Value = "'0.776"
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +3.41%  '
$c = $ws.Range('D45')
$c.Value = "'42.24"
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +3.93%  '
$ws.Range('E46').Value = '  +2.46%  '
$ws.Range('E47').Value = '  +10.17%  '
$ws.Range('E48').Value = '  +5.45%  '
$c = $ws.Range('D49')
$c.Value = "'2.23"
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +26.73%  '
$c = $ws.Range('D50')
$c.Value = "'0.864"
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +7.46%  '
$c = $ws.Range('D51')
$c.Value = "'6.54"
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +4.33%  '
